# daily auto push: 2026-01-08 13:50 UTC
#
# The source sheet is a flat log of timestamped events (date, weekday,
# hour, ranking). A new observation for 2026/01/08 (weekday 木, hour 20,
# ranking 201) was recorded between the existing "17:00" row and the
# "2026/12/29" row, so a single row is inserted at row 594 and every row
# that used to live at 594..635 shifts down to 595..636.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row above the old row 594 (2026/12/29 shifts to 595,
# ..., the old last row 635 shifts to 636). This also bumps the sheet's
# dimension from D635 to D636 automatically.
$ws.Rows("594:594").Insert()

# The date column in this workbook holds plain text like "2026/01/08",
# not a real Excel date. Typing that string into a General-formatted
# cell would normally be auto-parsed into a date serial number, so we
# briefly force text formatting while entering the value, then clear
# the formatting override again so the new cell ends up with no
# explicit style -- exactly like every other data row in the sheet.
$ws.Range("A594").NumberFormat = "@"
$ws.Range("A594").Value = "2026/01/08"
$ws.Range("A594").ClearFormats()

$ws.Range("B594").Value = "木"
$ws.Range("C594").Value = 20
$ws.Range("D594").Value = 201
